# Apply weekly update to fruit/vegetable price sheet (Pepino dulce)
# The data rows (2-14) get reshuffled: dates, volumes, quality and some
# prices move to different rows while the rest of the row content stays.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# P (Precio $/Kg)
$rows = @{
    2  = @{ D = 44645; I = "Primera"; J = 60;  K = 15000; L = 16000; M = 15500; P = 861 }
    3  = @{ D = 44630; I = "Primera"; J = 60;  K = 15000; L = 16000; M = 15500; P = 861 }
    4  = @{ D = 44384; I = "Segunda"; J = 60;  K = 15000; L = 15000; M = 15000; P = 833 }
    5  = @{ D = 44635; I = "Primera"; J = 100; K = 15000; L = 16000; M = 15500; P = 861 }
    6  = @{ D = 44642; I = "Primera"; J = 100; K = 15000; L = 16000; M = 15500; P = 861 }
    7  = @{ D = 44658; I = "Primera"; J = 80;  K = 15000; L = 16000; M = 15500; P = 861 }
    8  = @{ D = 44656; I = "Primera"; J = 100; K = 15000; L = 16000; M = 15500; P = 861 }
    9  = @{ D = 44637; I = "Primera"; J = 100; K = 15000; L = 16000; M = 15500; P = 861 }
    10 = @{ D = 44649; I = "Primera"; J = 60;  K = 15000; L = 16000; M = 15500; P = 861 }
    11 = @{ D = 44664; I = "Primera"; J = 160; K = 15000; L = 16000; M = 15500; P = 861 }
    12 = @{ D = 44628; I = "Primera"; J = 60;  K = 15000; L = 16000; M = 15500; P = 861 }
    13 = @{ D = 44659; I = "Primera"; J = 80;  K = 15000; L = 16000; M = 15500; P = 861 }
    14 = @{ D = 44651; I = "Primera"; J = 60;  K = 15000; L = 16000; M = 15500; P = 861 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 9).Value  = $vals.I   # I - Calidad
    $ws.Cells.Item($r, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio $/Kg
}
